# Simulated Wild Card round and logged it
# Update the "Road" (R) row of target-depth stats on both the OFF and DEF
# sheets to reflect the additional game that was played/logged.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 227
$wsOff.Range("C3").Value = 163
$wsOff.Range("D3").Value = 60
$wsOff.Range("E3").Value = 26

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 244
$wsDef.Range("C3").Value = 156
$wsDef.Range("D3").Value = 57
$wsDef.Range("E3").Value = 28
